$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Recommandations")
$ws2 = $wb.Worksheets.Item("Top_YTD")

# --- Sheet1 "Recommandations": update rows 2-37, delete old rows 38-40 ---
$ws1.Cells.Item(2, 1).Value = "BRVM - SERVICES PUBLICS"
$ws1.Cells.Item(2, 2).Value = 0
$ws1.Cells.Item(2, 3).Value = 6
$ws1.Cells.Item(2, 4).Value = 2543.62
$ws1.Cells.Item(2, 5).Value = 114.25
$ws1.Cells.Item(2, 6).Value = "🟡 Observer"
$ws1.Cells.Item(2, 7).Value = "➖ Neutre"
$ws1.Cells.Item(3, 1).Value = "NEI-CEDA CI"
$ws1.Cells.Item(3, 2).Value = 0
$ws1.Cells.Item(3, 3).Value = 3
$ws1.Cells.Item(3, 4).Value = 2390
$ws1.Cells.Item(3, 5).Value = 885
$ws1.Cells.Item(3, 6).Value = "🟡 Observer"
$ws1.Cells.Item(3, 7).Value = "➖ Neutre"
$ws1.Cells.Item(4, 1).Value = "AIR LIQUIDE CI"
$ws1.Cells.Item(4, 2).Value = 0
$ws1.Cells.Item(4, 3).Value = 3
$ws1.Cells.Item(4, 4).Value = 2105
$ws1.Cells.Item(4, 5).Value = 695
$ws1.Cells.Item(4, 6).Value = "🟡 Observer"
$ws1.Cells.Item(4, 7).Value = "➖ Neutre"
$ws1.Cells.Item(5, 1).Value = "BRVM - DISTRIBUTION"
$ws1.Cells.Item(5, 2).Value = 0
$ws1.Cells.Item(5, 3).Value = 3
$ws1.Cells.Item(5, 4).Value = 1815.82
$ws1.Cells.Item(5, 5).Value = 575.3
$ws1.Cells.Item(5, 6).Value = "🟡 Observer"
$ws1.Cells.Item(5, 7).Value = "➖ Neutre"
$ws1.Cells.Item(6, 1).Value = "BRVM - AUTRES SECTEURS"
$ws1.Cells.Item(6, 2).Value = 0
$ws1.Cells.Item(6, 3).Value = 3
$ws1.Cells.Item(6, 4).Value = 1803.92
$ws1.Cells.Item(6, 5).Value = 597.46
$ws1.Cells.Item(6, 6).Value = "🟡 Observer"
$ws1.Cells.Item(6, 7).Value = "➖ Neutre"
$ws1.Cells.Item(7, 1).Value = "BRVM - TRANSPORT"
$ws1.Cells.Item(7, 2).Value = 0
$ws1.Cells.Item(7, 3).Value = 3
$ws1.Cells.Item(7, 4).Value = 1079.45
$ws1.Cells.Item(7, 5).Value = 364.71
$ws1.Cells.Item(7, 6).Value = "🟡 Observer"
$ws1.Cells.Item(7, 7).Value = "➖ Neutre"
$ws1.Cells.Item(8, 1).Value = "BRVM - AGRICULTURE"
$ws1.Cells.Item(8, 2).Value = 0
$ws1.Cells.Item(8, 3).Value = 3
$ws1.Cells.Item(8, 4).Value = 1026.56
$ws1.Cells.Item(8, 5).Value = 332.84
$ws1.Cells.Item(8, 6).Value = "🟡 Observer"
$ws1.Cells.Item(8, 7).Value = "➖ Neutre"
$ws1.Cells.Item(9, 1).Value = "BRVM - CONSOMMATION DISCRETIONNAIRE"
$ws1.Cells.Item(9, 2).Value = 0
$ws1.Cells.Item(9, 3).Value = 3
$ws1.Cells.Item(9, 4).Value = 689.62
$ws1.Cells.Item(9, 5).Value = 215.49
$ws1.Cells.Item(9, 6).Value = "🟡 Observer"
$ws1.Cells.Item(9, 7).Value = "➖ Neutre"
$ws1.Cells.Item(10, 1).Value = "BRVM-PRESTIGE"
$ws1.Cells.Item(10, 2).Value = 0
$ws1.Cells.Item(10, 3).Value = 3
$ws1.Cells.Item(10, 4).Value = 430.05
$ws1.Cells.Item(10, 5).Value = 141.58
$ws1.Cells.Item(10, 6).Value = "🟡 Observer"
$ws1.Cells.Item(10, 7).Value = "➖ Neutre"
$ws1.Cells.Item(11, 1).Value = "BRVM - FINANCES"
$ws1.Cells.Item(11, 2).Value = 0
$ws1.Cells.Item(11, 3).Value = 3
$ws1.Cells.Item(11, 4).Value = 429.34
$ws1.Cells.Item(11, 5).Value = 141.7
$ws1.Cells.Item(11, 6).Value = "🟡 Observer"
$ws1.Cells.Item(11, 7).Value = "➖ Neutre"
$ws1.Cells.Item(12, 1).Value = "BRVM - SERVICES FINANCIERS"
$ws1.Cells.Item(12, 2).Value = 0
$ws1.Cells.Item(12, 3).Value = 3
$ws1.Cells.Item(12, 4).Value = 421.95
$ws1.Cells.Item(12, 5).Value = 139.26
$ws1.Cells.Item(12, 6).Value = "🟡 Observer"
$ws1.Cells.Item(12, 7).Value = "➖ Neutre"
$ws1.Cells.Item(13, 1).Value = "BRVM - INDUSTRIELS"
$ws1.Cells.Item(13, 2).Value = 0
$ws1.Cells.Item(13, 3).Value = 3
$ws1.Cells.Item(13, 4).Value = 370.37
$ws1.Cells.Item(13, 5).Value = 122.24
$ws1.Cells.Item(13, 6).Value = "🟡 Observer"
$ws1.Cells.Item(13, 7).Value = "➖ Neutre"
$ws1.Cells.Item(14, 1).Value = "BRVM - ENERGIE"
$ws1.Cells.Item(14, 2).Value = 0
$ws1.Cells.Item(14, 3).Value = 3
$ws1.Cells.Item(14, 4).Value = 321.22
$ws1.Cells.Item(14, 5).Value = 105.97
$ws1.Cells.Item(14, 6).Value = "🟡 Observer"
$ws1.Cells.Item(14, 7).Value = "➖ Neutre"
$ws1.Cells.Item(15, 1).Value = "BRVM - TELECOMMUNICATIONS"
$ws1.Cells.Item(15, 2).Value = 0
$ws1.Cells.Item(15, 3).Value = 3
$ws1.Cells.Item(15, 4).Value = 286.33
$ws1.Cells.Item(15, 5).Value = 95.22
$ws1.Cells.Item(15, 6).Value = "🟡 Observer"
$ws1.Cells.Item(15, 7).Value = "➖ Neutre"
$ws1.Cells.Item(16, 1).Value = "BRVM - INDUSTRIE                (**)"
$ws1.Cells.Item(16, 2).Value = 0
$ws1.Cells.Item(16, 3).Value = 1
$ws1.Cells.Item(16, 4).Value = 235.98
$ws1.Cells.Item(16, 5).Value = 235.98
$ws1.Cells.Item(16, 6).Value = "🟡 Observer"
$ws1.Cells.Item(16, 7).Value = "➖ Neutre"
$ws1.Cells.Item(17, 1).Value = "BRVM-PRINCIPAL                    (**)"
$ws1.Cells.Item(17, 2).Value = 0
$ws1.Cells.Item(17, 3).Value = 1
$ws1.Cells.Item(17, 4).Value = 213.74
$ws1.Cells.Item(17, 5).Value = 213.74
$ws1.Cells.Item(17, 6).Value = "🟡 Observer"
$ws1.Cells.Item(17, 7).Value = "➖ Neutre"
$ws1.Cells.Item(18, 1).Value = "BRVM - CONSOMMATION DE BASE         (**)"
$ws1.Cells.Item(18, 2).Value = 0
$ws1.Cells.Item(18, 3).Value = 1
$ws1.Cells.Item(18, 4).Value = 202.47
$ws1.Cells.Item(18, 5).Value = 202.47
$ws1.Cells.Item(18, 6).Value = "🟡 Observer"
$ws1.Cells.Item(18, 7).Value = "➖ Neutre"
$ws1.Cells.Item(19, 1).Value = "SOLIBRA CI (SLBC)"
$ws1.Cells.Item(19, 2).Value = 3
$ws1.Cells.Item(19, 3).Value = 0
$ws1.Cells.Item(19, 4).Value = 22.47
$ws1.Cells.Item(19, 5).Value = 7.49
$ws1.Cells.Item(19, 6).Value = "🟢 Achat"
$ws1.Cells.Item(19, 7).Value = "✅ Renforcer"
$ws1.Cells.Item(20, 1).Value = "NEI-CEDA CI (NEIC)"
$ws1.Cells.Item(20, 2).Value = 3
$ws1.Cells.Item(20, 3).Value = 0
$ws1.Cells.Item(20, 4).Value = 21.35
$ws1.Cells.Item(20, 5).Value = 7.27
$ws1.Cells.Item(20, 6).Value = "🟢 Achat"
$ws1.Cells.Item(20, 7).Value = "✅ Renforcer"
$ws1.Cells.Item(21, 1).Value = "SICABLE CI (CABC)"
$ws1.Cells.Item(21, 2).Value = 2
$ws1.Cells.Item(21, 3).Value = 0
$ws1.Cells.Item(21, 4).Value = 14.55
$ws1.Cells.Item(21, 5).Value = 7.3
$ws1.Cells.Item(21, 6).Value = "🟡 Observer"
$ws1.Cells.Item(21, 7).Value = "➖ Neutre"
$ws1.Cells.Item(22, 1).Value = "CFAO MOTORS CI (CFAC)"
$ws1.Cells.Item(22, 2).Value = 1
$ws1.Cells.Item(22, 3).Value = 0
$ws1.Cells.Item(22, 4).Value = 7.08
$ws1.Cells.Item(22, 5).Value = 7.08
$ws1.Cells.Item(22, 6).Value = "🟡 Observer"
$ws1.Cells.Item(22, 7).Value = "➖ Neutre"
$ws1.Cells.Item(23, 1).Value = "UNILEVER CI (UNLC)"
$ws1.Cells.Item(23, 2).Value = 1
$ws1.Cells.Item(23, 3).Value = 0
$ws1.Cells.Item(23, 4).Value = 6.98
$ws1.Cells.Item(23, 5).Value = 6.98
$ws1.Cells.Item(23, 6).Value = "🟡 Observer"
$ws1.Cells.Item(23, 7).Value = "➖ Neutre"
$ws1.Cells.Item(24, 1).Value = "SOCIETE IVOIRIENNE DE BANQUE  (SIBC)"
$ws1.Cells.Item(24, 2).Value = 1
$ws1.Cells.Item(24, 3).Value = 0
$ws1.Cells.Item(24, 4).Value = 2.59
$ws1.Cells.Item(24, 5).Value = 2.59
$ws1.Cells.Item(24, 6).Value = "🟡 Observer"
$ws1.Cells.Item(24, 7).Value = "➖ Neutre"
$ws1.Cells.Item(25, 1).Value = "SITAB CI (STBC)"
$ws1.Cells.Item(25, 2).Value = 1
$ws1.Cells.Item(25, 3).Value = 0
$ws1.Cells.Item(25, 4).Value = 2.56
$ws1.Cells.Item(25, 5).Value = 2.56
$ws1.Cells.Item(25, 6).Value = "🟡 Observer"
$ws1.Cells.Item(25, 7).Value = "➖ Neutre"
$ws1.Cells.Item(26, 1).Value = "AFRICA GLOBAL LOGISTICS CI (SDSC)"
$ws1.Cells.Item(26, 2).Value = 1
$ws1.Cells.Item(26, 3).Value = 0
$ws1.Cells.Item(26, 4).Value = 1.36
$ws1.Cells.Item(26, 5).Value = 1.36
$ws1.Cells.Item(26, 6).Value = "🟡 Observer"
$ws1.Cells.Item(26, 7).Value = "➖ Neutre"
$ws1.Cells.Item(27, 1).Value = "TOTAL"
$ws1.Cells.Item(27, 2).Value = 0
$ws1.Cells.Item(27, 3).Value = 3
$ws1.Cells.Item(27, 4).Value = 0
$ws1.Cells.Item(27, 5).Value = 0
$ws1.Cells.Item(27, 6).Value = "🟡 Observer"
$ws1.Cells.Item(27, 7).Value = "➖ Neutre"
$ws1.Cells.Item(28, 1).Value = "SAPH CI (SPHC)"
$ws1.Cells.Item(28, 2).Value = 1
$ws1.Cells.Item(28, 3).Value = 1
$ws1.Cells.Item(28, 4).Value = -2
$ws1.Cells.Item(28, 5).Value = 2.53
$ws1.Cells.Item(28, 6).Value = "🟡 Observer"
$ws1.Cells.Item(28, 7).Value = "👀 À surveiller"
$ws1.Cells.Item(29, 1).Value = "SOGB CI (SOGC)"
$ws1.Cells.Item(29, 2).Value = 0
$ws1.Cells.Item(29, 3).Value = 1
$ws1.Cells.Item(29, 4).Value = -3.83
$ws1.Cells.Item(29, 5).Value = -3.83
$ws1.Cells.Item(29, 6).Value = "🟡 Observer"
$ws1.Cells.Item(29, 7).Value = "➖ Neutre"
$ws1.Cells.Item(30, 1).Value = "LOTERIE NATIONALE DU BENIN (LNBB)"
$ws1.Cells.Item(30, 2).Value = 0
$ws1.Cells.Item(30, 3).Value = 1
$ws1.Cells.Item(30, 4).Value = -4.63
$ws1.Cells.Item(30, 5).Value = -4.63
$ws1.Cells.Item(30, 6).Value = "🟡 Observer"
$ws1.Cells.Item(30, 7).Value = "➖ Neutre"
$ws1.Cells.Item(31, 1).Value = "TRACTAFRIC MOTORS CI (PRSC)"
$ws1.Cells.Item(31, 2).Value = 1
$ws1.Cells.Item(31, 3).Value = 2
$ws1.Cells.Item(31, 4).Value = -7.43
$ws1.Cells.Item(31, 5).Value = -7.42
$ws1.Cells.Item(31, 6).Value = "🟡 Observer"
$ws1.Cells.Item(31, 7).Value = "👀 À surveiller"
$ws1.Cells.Item(32, 1).Value = "BERNABE CI (BNBC)"
$ws1.Cells.Item(32, 2).Value = 0
$ws1.Cells.Item(32, 3).Value = 1
$ws1.Cells.Item(32, 4).Value = -7.44
$ws1.Cells.Item(32, 5).Value = -7.44
$ws1.Cells.Item(32, 6).Value = "🟡 Observer"
$ws1.Cells.Item(32, 7).Value = "➖ Neutre"
$ws1.Cells.Item(33, 1).Value = "BANK OF AFRICA ML (BOAM)"
$ws1.Cells.Item(33, 2).Value = 0
$ws1.Cells.Item(33, 3).Value = 1
$ws1.Cells.Item(33, 4).Value = -7.47
$ws1.Cells.Item(33, 5).Value = -7.47
$ws1.Cells.Item(33, 6).Value = "🟡 Observer"
$ws1.Cells.Item(33, 7).Value = "➖ Neutre"
$ws1.Cells.Item(34, 1).Value = "SAFCA CI (SAFC)"
$ws1.Cells.Item(34, 2).Value = 0
$ws1.Cells.Item(34, 3).Value = 1
$ws1.Cells.Item(34, 4).Value = -7.5
$ws1.Cells.Item(34, 5).Value = -7.5
$ws1.Cells.Item(34, 6).Value = "🟡 Observer"
$ws1.Cells.Item(34, 7).Value = "➖ Neutre"
$ws1.Cells.Item(35, 1).Value = "PALM CI (PALC)"
$ws1.Cells.Item(35, 2).Value = 0
$ws1.Cells.Item(35, 3).Value = 2
$ws1.Cells.Item(35, 4).Value = -14.94
$ws1.Cells.Item(35, 5).Value = -7.49
$ws1.Cells.Item(35, 6).Value = "🟡 Observer"
$ws1.Cells.Item(35, 7).Value = "➖ Neutre"
$ws1.Cells.Item(36, 1).Value = "NESTLE CI (NTLC)"
$ws1.Cells.Item(36, 2).Value = 0
$ws1.Cells.Item(36, 3).Value = 2
$ws1.Cells.Item(36, 4).Value = -14.96
$ws1.Cells.Item(36, 5).Value = -7.48
$ws1.Cells.Item(36, 6).Value = "🟡 Observer"
$ws1.Cells.Item(36, 7).Value = "➖ Neutre"
$ws1.Cells.Item(37, 1).Value = "FILTISAC CI (FTSC)"
$ws1.Cells.Item(37, 2).Value = 0
$ws1.Cells.Item(37, 3).Value = 3
$ws1.Cells.Item(37, 4).Value = -22.12
$ws1.Cells.Item(37, 5).Value = -7.48
$ws1.Cells.Item(37, 6).Value = "🔴 Vente"
$ws1.Cells.Item(37, 7).Value = "⚠️ Risque de décrochage"

# Remove now-unused trailing rows (38-40), shrinking dimension to A1:G37
$ws1.Range("A38:A40").EntireRow.Delete()

# --- Sheet2 "Top_YTD": update progression values ---
$ws2.Cells.Item(2, 2).Value = 568629.85
$ws2.Cells.Item(3, 2).Value = 71455.32000000001
$ws2.Cells.Item(4, 2).Value = 51416
$ws2.Cells.Item(5, 2).Value = 34920.46
$ws2.Cells.Item(6, 2).Value = 34391.37
$ws2.Cells.Item(7, 2).Value = 9620.860000000001
$ws2.Cells.Item(8, 2).Value = 8543.1
$ws2.Cells.Item(9, 2).Value = 3482.69
$ws2.Cells.Item(10, 2).Value = 1341.03
$ws2.Cells.Item(11, 2).Value = 1336.86
